$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New master data rows (regcntr_id, usr_id, machine_id, lang_code, is_active, cr_by, cr_dtimes)
$newRows = @(
    @(10005, 110033, 10005, "eng", $true, "superadmin", "now()"),
    @(10005, 110034, 10005, "eng", $true, "superadmin", "now()"),
    @(10005, 110035, 10005, "eng", $true, "superadmin", "now()")
)

$startRow = 34
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}

# Move selection to the row after the newly-added data, matching the
# post-edit state (full-row selection on the next empty row).
$ws.Range("A37").Select()
$ws.Range("A37:XFD1048576").Select()
